$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number for data rows 2-91.
# The diff changes every one of these values from 45182 (2023-09-13) to
# 45184 (2023-09-15).
for ($row = 2; $row -le 91; $row++) {
    $ws.Cells.Item($row, 3).Value = 45184
}
